$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in previously empty cells in row 10 with placeholder values
$ws.Range("D10").Value = "NA"
$ws.Range("F10").Value = "placeholder"

# Update the saved selection/active cell for the sheet
$ws.Range("F18").Select()
